$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.066.60'
$ws.Range("E2").Value = '  -0.98%  '
$ws.Range("D3").Value = '1.790.52'
$ws.Range("E3").Value = '  -1.50%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = "'227.02"
$ws.Range("E5").Value = '  -1.58%  '
$ws.Range("E6").Value = '  +1.55%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").Value = "'31.26"
$ws.Range("E8").Value = '  -0.72%  '
$ws.Range("D9").Value = "'46.03"
$ws.Range("E9").Value = '  +0.53%  '
$ws.Range("E10").Value = '  -0.64%  '
$ws.Range("E11").Value = '  -2.97%  '
$ws.Range("D12").Value = "'0.0927"
$ws.Range("E12").Value = '  -0.47%  '
$ws.Range("D13").Value = '2.047.11'
$ws.Range("E13").Value = '  -1.60%  '
$ws.Range("E14").Value = '  +11.97%  '
$ws.Range("D15").Value = '1.782.49'
$ws.Range("E15").Value = '  -2.16%  '
$ws.Range("D16").Value = "'0.637"
$ws.Range("E16").Value = '  -1.30%  '
$ws.Range("D17").Value = '34.068.26'
$ws.Range("E17").Value = '  -0.77%  '
$ws.Range("D18").Value = "'4.23"
$ws.Range("E18").Value = '  -2.67%  '
$ws.Range("D19").Value = "'69.61"
$ws.Range("E19").Value = '  -1.22%  '
$ws.Range("D20").Value = "'253.14"
$ws.Range("E20").Value = '  -3.02%  '
$ws.Range("D21").Value = '0.0₃0743'
$ws.Range("E21").Value = '  -1.26%  '
$ws.Range("E23").Value = '  -0.64%  '
$ws.Range("E24").Value = '  -2.02%  '
$ws.Range("E25").Value = '  -2.43%  '
$ws.Range("D26").Value = "'156.97"
$ws.Range("D27").Value = "'16.60"
$ws.Range("E27").Value = '  -1.24%  '
$ws.Range("E28").Value = '  -1.62%  '
$ws.Range("E29").Value = '  -2.12%  '
$ws.Range("E30").Value = '  -0.07%  '
$ws.Range("D31").Value = "'3.82"
$ws.Range("E31").Value = '  -0.08%  '
$ws.Range("E32").Value = '  +0.05%  '
$ws.Range("E33").Value = '  -1.34%  '
$ws.Range("D34").Value = "'3.61"
$ws.Range("E34").Value = '  +1.18%  '
$ws.Range("D35").Value = "'1.85"
$ws.Range("E35").Value = '  +0.57%  '
$ws.Range("D36").Value = '1.450.16'
$ws.Range("E36").Value = '  -8.45%  '
$ws.Range("E37").Value = '  +0.74%  '
$ws.Range("D38").Value = "'0.633"
$ws.Range("E38").Value = '  +0.08%  '
$ws.Range("E39").Value = '  -1.31%  '
$ws.Range("D40").Value = "'83.22"
$ws.Range("E40").Value = '  -2.59%  '
$ws.Range("D41").Value = "'2.82"
$ws.Range("E41").Value = '  -0.65%  '
$ws.Range("E42").Value = '  +0.18%  '
$ws.Range("E43").Value = '  -1.56%  '
$ws.Range("E44").Value = '  -2.42%  '
$ws.Range("E45").Value = '  -1.89%  '
$ws.Range("E46").Value = '  +0.66%  '
$ws.Range("D47").Value = '1.948.14'
$ws.Range("E47").Value = '  -1.63%  '
$ws.Range("E48").Value = '  -0.19%  '
$ws.Range("E49").Value = '  +0.08%  '
$ws.Range("D50").Value = "'11.83"
$ws.Range("E50").Value = '  +4.97%  '
$ws.Range("D51").Value = "'51.27"
$ws.Range("E51").Value = '  -4.13%  '
